# Bugfixed the naive forecaster component module
#
# The original first data row (dated 2007-10-... with no forecast values)
# was erroneous and needs to be dropped; every remaining data row moves up
# by one, and the y_1_forecast (column E) figures are recomputed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete first data row (row 2); this shifts rows 3:19 up to 2:18
# and automatically keeps A/B/C/D values intact while updating the sheet
# dimension.
$ws.Range("A2").EntireRow.Delete()

# Recompute the y_1_forecast (column E) values for the now-shifted rows.
$newE = @(
    0.4282194198276246,
    -0.8235211753995442,
    1.121293995080253,
    1.665250327443002,
    1.079796209653616,
    1.374377011838535,
    1.310895847186577,
    1.862478303083726,
    1.639776099317536,
    2.181728312936415,
    2.010025322622599,
    0.7771393814490102,
    -2.71887004062904,
    0.4582698374457683,
    1.314675624401973,
    0.1199358335146838,
    0.2676745853112728
)

for ($i = 0; $i -lt $newE.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $newE[$i]
}
